# Adds 'drop' and 'dropExplanation' columns (Q, R) to Sheet2, mirroring the
# header styling used by the existing J1:P1 header cells but with a
# slightly different fill/border, and makes Sheet2 the active tab.

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Header cells (row 1) -------------------------------------------------
$q1 = $ws2.Range("Q1")
$r1 = $ws2.Range("R1")

# Start from the existing header style (font + fill + border + alignment),
# then tweak the fill's background colour and the border sides so the
# result matches the new header look. Do this BEFORE writing the header
# text, since Range.Copy also copies the source cell's value.
$k1 = $ws2.Range("K1")
[void]$k1.Copy($q1)

$q1.Interior.Color = 16764057
$q1.Interior.PatternColor = 0

[void]$q1.Copy($r1)
$r1.Borders.Item(7).LineStyle = -4142
$r1.Borders.Item(10).LineStyle = 1
$r1.Borders.Item(10).Color = 0

$q1.Value = "drop"
$r1.Value = "dropExplanation"

# --- Data cells (rows 2-41) ------------------------------------------------
# Q holds a boolean "drop" flag (defaults to FALSE); R is the (blank)
# explanation column. Both just need the plain-black-text font tweak.
$dataQ = $ws2.Range("Q2:Q41")
$dataQ.Value = $false

$dataBoth = $ws2.Range("Q2:R41")
$dataBoth.Font.Color = 0

# --- Active-tab bookkeeping -------------------------------------------------
# The author was on Sheet2 (with Q1:R41 selected) when the workbook was saved.
[void]$ws2.Activate()
[void]$ws2.Range("Q1:R41").Select()
